$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1381.3
$ws.Range("I58").Value = 238
$ws.Range("J58").Value = 1871.2858
$ws.Range("K58").Value = 714
$ws.Range("L58").Value = 5613.857400000001
$ws.Range("M58").Value = -564
$ws.Range("N58").Value = -5913.857400000001
$ws.Range("H86").Value = 4221.25
$ws.Range("I86").Value = 3193
$ws.Range("K86").Value = 3193
$ws.Range("M86").Value = -2070
$ws.Range("H89").Value = 4221.25
$ws.Range("I89").Value = 3193
$ws.Range("K89").Value = 15965
$ws.Range("M89").Value = -10349
$ws.Range("H138").Value = 5845.8887
$ws.Range("I138").Value = 1999
$ws.Range("J138").Value = 6326.75
$ws.Range("K138").Value = 5997
$ws.Range("L138").Value = 18980.25
$ws.Range("M138").Value = -857
$ws.Range("N138").Value = -29260.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2468.375
$ws.Range("I132").Value = 2487.2
$ws.Range("K132").Value = 7461.599999999999
$ws.Range("M132").Value = -4931.599999999999
$ws.Range("H134").Value = 79999.664
$ws.Range("I134").Value = 70000
$ws.Range("K134").Value = 70000
$ws.Range("M134").Value = -64930
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 607.0714
$ws.Range("I80").Value = 357.7
$ws.Range("K80").Value = 357.7
$ws.Range("M80").Value = 640.3
$ws.Range("H81").Value = 77333.336
$ws.Range("J81").Value = 77333.336
$ws.Range("L81").Value = 77333.336
$ws.Range("N81").Value = -79455.336
$ws.Range("H83").Value = 607.0714
$ws.Range("I83").Value = 357.7
$ws.Range("K83").Value = 1788.5
$ws.Range("M83").Value = 3203.5
$ws.Range("H84").Value = 77333.336
$ws.Range("J84").Value = 77333.336
$ws.Range("L84").Value = 232000.008
$ws.Range("N84").Value = -242608.008
$ws.Range("H108").Value = 37522
$ws.Range("J108").Value = 37522
$ws.Range("L108").Value = 37522
$ws.Range("N108").Value = -45202
$ws.Range("H134").Value = 3867.6924
$ws.Range("I134").Value = 3844.5454
$ws.Range("K134").Value = 11533.6362
$ws.Range("M134").Value = -8998.636200000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 116.666664
$ws.Range("I6").Value = 116.666664
$ws.Range("K6").Value = 116.666664
$ws.Range("M6").Value = -3.666663999999997
$ws.Range("H9").Value = 51401.332
$ws.Range("J9").Value = 51401.332
$ws.Range("L9").Value = 51401.332
$ws.Range("N9").Value = -51737.332
$ws.Range("H22").Value = 102
$ws.Range("I22").Value = 102
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 102
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 248
$ws.Range("N22").ClearContents()
$ws.Range("H31").Value = 3888.05
$ws.Range("J31").Value = 4802.7383
$ws.Range("L31").Value = 4802.7383
$ws.Range("N31").Value = -5392.7383
$ws.Range("H34").Value = 3888.05
$ws.Range("J34").Value = 4802.7383
$ws.Range("L34").Value = 4802.7383
$ws.Range("N34").Value = -5206.7383
$ws.Range("H94").Value = 2957.5334
$ws.Range("I94").Value = 1020
$ws.Range("J94").Value = 6832.6
$ws.Range("K94").Value = 1020
$ws.Range("L94").Value = 6832.6
$ws.Range("M94").Value = -569
$ws.Range("N94").Value = -7734.6
$ws.Range("H132").Value = 1766.6666
$ws.Range("I132").Value = 1766.6666
$ws.Range("K132").Value = 5299.9998
$ws.Range("M132").Value = -2769.9998
$ws.Range("H134").Value = 461.69232
$ws.Range("I134").Value = 461.69232
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 1385.07696
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 1149.92304
$ws.Range("N134").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 61.6
$ws.Range("I2").Value = 28.5
$ws.Range("J2").Value = 127.8
$ws.Range("K2").Value = 171
$ws.Range("L2").Value = 766.8
$ws.Range("M2").Value = -58
$ws.Range("N2").Value = -992.8
$ws.Range("H38").Value = 344.70834
$ws.Range("I38").Value = 391.78946
$ws.Range("J38").Value = 165.8
$ws.Range("K38").Value = 1175.36838
$ws.Range("L38").Value = 497.4
$ws.Range("M38").Value = -828.3683800000001
$ws.Range("N38").Value = -1191.4
$ws.Range("H40").Value = 338.5
$ws.Range("J40").Value = 338.5
$ws.Range("L40").Value = 1354
$ws.Range("N40").Value = -1492
$ws.Range("H141").Value = 846
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 15853625
$ws.Range("I3").Value = 21777988
$ws.Range("J3").Value = 4004900
$ws.Range("K3").Value = 21777988
$ws.Range("L3").Value = 4004900
$ws.Range("M3").Value = -21777872
$ws.Range("N3").Value = -4005132
$ws.Range("H132").Value = 3263.4614
$ws.Range("I132").Value = 3263.4614
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9790.3842
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7260.3842
$ws.Range("N132").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 5000
$ws.Range("K3").Value = 5000
$ws.Range("M3").Value = -4888
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("I14").Value = 13443.444
$ws.Range("J14").Value = 13199.6
$ws.Range("K14").Value = 13443.444
$ws.Range("L14").Value = 13199.6
$ws.Range("M14").Value = -13271.444
$ws.Range("N14").Value = -13543.6
$ws.Range("H15").Value = 5000
$ws.Range("I15").Value = 5000
$ws.Range("K15").Value = 5000
$ws.Range("M15").Value = -4830
$ws.Range("H82").Value = 5638.1113
$ws.Range("I82").Value = 3936
$ws.Range("J82").Value = 6999.8
$ws.Range("K82").Value = 3936
$ws.Range("L82").Value = 6999.8
$ws.Range("M82").Value = -3575
$ws.Range("N82").Value = -7721.8
$ws.Range("H85").Value = 5638.1113
$ws.Range("I85").Value = 3936
$ws.Range("J85").Value = 6999.8
$ws.Range("K85").Value = 3936
$ws.Range("L85").Value = 6999.8
$ws.Range("M85").Value = -2688
$ws.Range("N85").Value = -9495.799999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 9997.5
$ws.Range("J13").Value = 9997.5
$ws.Range("L13").Value = 9997.5
$ws.Range("N13").Value = -10277.5
$ws.Range("H52").Value = 26663
$ws.Range("J52").Value = 20000
$ws.Range("L52").Value = 20000
$ws.Range("N52").Value = -20452
$ws.Range("H125").Value = 17633
$ws.Range("J125").Value = 17633
$ws.Range("L125").Value = 17633
$ws.Range("N125").Value = -27473
$ws.Range("H136").Value = 2083.9033
$ws.Range("I136").Value = 1102.4762
$ws.Range("K136").Value = 3307.4286
$ws.Range("M136").Value = -757.4286000000002
$ws.Range("H137").Value = 61665
$ws.Range("J137").Value = 61665
$ws.Range("L137").Value = 61665
$ws.Range("N137").Value = -71865
